$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3014
$ws.Range("F3").Value = 3014
$ws.Range("F4").Value = 6436
$ws.Range("F5").Value = 2559
$ws.Range("F6").Value = 658
$ws.Range("F7").Value = 89
$ws.Range("F8").Value = 3151
$ws.Range("F9").Value = 371
$ws.Range("F10").Value = 43
$ws.Range("F11").Value = 7734
$ws.Range("F12").Value = 389
$ws.Range("F13").Value = 70
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 7
$ws.Range("F17").Value = 264
$ws.Range("F19").Value = 25
$ws.Range("F20").Value = 487
$ws.Range("F21").Value = 9525
$ws.Range("F22").Value = 26
$ws.Range("F28").Value = 133
$ws.Range("F31").Value = 129
$ws.Range("F32").Value = 75
$ws.Range("F33").Value = 2628
$ws.Range("F37").Value = 1493
$ws.Range("F38").Value = 807
$ws.Range("F39").Value = 3979
$ws.Range("F41").Value = 696
$ws.Range("F43").Value = 120
$ws.Range("F44").Value = 260
$ws.Range("F45").Value = 88
$ws.Range("F46").Value = 19
$ws.Range("F48").Value = 44
$ws.Range("F50").Value = 19

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 23
$ws.Range("F19").Value = 169
$ws.Range("F21").Value = 16

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3014
$ws.Range("F3").Value = 3014
$ws.Range("F6").Value = 6436
$ws.Range("F7").Value = 2559
$ws.Range("F8").Value = 658
$ws.Range("F9").Value = 89
$ws.Range("F10").Value = 3151
$ws.Range("F11").Value = 371
$ws.Range("F14").Value = 43
$ws.Range("F15").Value = 7734
$ws.Range("F16").Value = 389
$ws.Range("F17").Value = 70
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 264
$ws.Range("F21").Value = 25
$ws.Range("F22").Value = 9525
$ws.Range("F23").Value = 26
$ws.Range("F27").Value = 133
$ws.Range("F29").Value = 129
$ws.Range("F30").Value = 75
$ws.Range("F31").Value = 2628
$ws.Range("F34").Value = 1493
$ws.Range("F35").Value = 807
$ws.Range("F36").Value = 169
$ws.Range("F37").Value = 3979
$ws.Range("F39").Value = 697
$ws.Range("F40").Value = 16
$ws.Range("F42").Value = 120
$ws.Range("F43").Value = 260
$ws.Range("F45").Value = 88
$ws.Range("F46").Value = 19
$ws.Range("F48").Value = 44
$ws.Range("F50").Value = 19
